$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '27.227.13'
Set-TextValue "E2" '  -2.75%  '

Set-TextValue "D3" '1.702.65'
Set-TextValue "E3" '  -2.14%  '

Set-TextValue "D4" '1.002'
Set-TextValue "E4" '  +0.00%  '

Set-TextValue "D5" '222.98'
Set-TextValue "E5" '  -2.69%  '

Set-TextValue "D6" '0.5290'
Set-TextValue "E6" '  -3.31%  '

Set-TextValue "D7" '1.002'
Set-TextValue "E7" '  +0.00%  '

Set-TextValue "D8" '0.2644'
Set-TextValue "E8" '  -5.13%  '

Set-TextValue "D9" '0.06572'
Set-TextValue "E9" '  -2.36%  '

Set-TextValue "D10" '20.77'
Set-TextValue "E10" '  -4.89%  '

Set-TextValue "D11" '0.07617'
Set-TextValue "E11" '  -2.32%  '

Set-TextValue "D12" '4.562'
Set-TextValue "E12" '  -3.14%  '

Set-TextValue "B13" 'WrappedEther'
Set-TextValue "C13" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D13" '1.695.09'
Set-TextValue "E13" '  -1.80%  '

Set-TextValue "B14" 'WrappedliquidstakedEther2.0'
Set-TextValue "C14" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D14" '1.937.74'
Set-TextValue "E14" '  -2.06%  '

Set-TextValue "D15" '0.5711'
Set-TextValue "E15" '  -5.21%  '

Set-TextValue "D16" '0.0₅8145'
Set-TextValue "E16" '  -3.42%  '

Set-TextValue "D17" '67.44'
Set-TextValue "E17" '  -3.21%  '

Set-TextValue "D18" '27.203.29'
Set-TextValue "E18" '  -2.75%  '

Set-TextValue "D19" '214.82'
Set-TextValue "E19" '  -5.25%  '

Set-TextValue "D20" '1.002'
Set-TextValue "E20" '  -0.07%  '

Set-TextValue "D21" '4.653'
Set-TextValue "E21" '  -4.09%  '

Set-TextValue "D22" '10.39'
Set-TextValue "E22" '  -5.42%  '

Set-TextValue "D23" '5.945'
Set-TextValue "E23" '  -4.86%  '

Set-TextValue "D24" '1.002'
Set-TextValue "E24" '  -0.05%  '

Set-TextValue "D25" '141.30'
Set-TextValue "E25" '  -3.54%  '

Set-TextValue "D26" '1.753'
Set-TextValue "E26" '  +5.91%  '

Set-TextValue "D27" '0.1212'
Set-TextValue "E27" '  -3.36%  '

Set-TextValue "D28" '7.234'
Set-TextValue "E28" '  -3.28%  '

Set-TextValue "D29" '16.26'
Set-TextValue "E29" '  -5.26%  '

Set-TextValue "D30" '0.05365'
Set-TextValue "E30" '  -5.66%  '

Set-TextValue "D31" '1.288'
Set-TextValue "E31" '  -2.29%  '

Set-TextValue "D32" '3.484'
Set-TextValue "E32" '  -5.98%  '

Set-TextValue "D33" '3.407'
Set-TextValue "E33" '  -4.01%  '

Set-TextValue "D34" '1.636'
Set-TextValue "E34" '  -1.99%  '

Set-TextValue "D35" '2.867'
Set-TextValue "E35" '  +0.21%  '

Set-TextValue "D36" '2.421'
Set-TextValue "E36" '  -1.18%  '

Set-TextValue "D37" '0.9443'
Set-TextValue "E37" '  -4.20%  '

Set-TextValue "D38" '0.5826'
Set-TextValue "E38" '  -2.49%  '

Set-TextValue "D39" '0.01625'
Set-TextValue "E39" '  -3.14%  '

Set-TextValue "D40" '5.859'
Set-TextValue "E40" '  -2.61%  '

Set-TextValue "D41" '1.002'
Set-TextValue "E41" '  +0.00%  '

Set-TextValue "D42" '1.040.50'
Set-TextValue "E42" '  -0.88%  '

Set-TextValue "D43" '0.8354'
Set-TextValue "E43" '  -1.32%  '

Set-TextValue "D44" '100.63'
Set-TextValue "E44" '  -1.60%  '

Set-TextValue "D45" '1.844.27'
Set-TextValue "E45" '  -2.10%  '

Set-TextValue "D46" '0.0₈114'
Set-TextValue "E46" '  -2.47%  '

Set-TextValue "D47" '57.85'
Set-TextValue "E47" '  -4.02%  '

Set-TextValue "D48" '0.4496'
Set-TextValue "E48" '  +1.54%  '

Set-TextValue "D49" '1.002'
Set-TextValue "E49" '  -0.81%  '

Set-TextValue "D50" '8.078'
Set-TextValue "E50" '  -2.96%  '

Set-TextValue "D51" '0.05236'
Set-TextValue "E51" '  -1.53%  '
